$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r0")
$ws.Activate()

$ws.Range("A23").Value = "ExoT_r0_script_9v1"
$ws.Range("B23").Value = "Andrew Lin 081120"
$ws.Range("C23").Value = "5 mL"
$ws.Range("D23").Value = "5 mL"
$ws.Range("E23").Value = "1 hour"
$ws.Range("F23").Value = "0.5 mL"
$ws.Range("G23").Value = "5 mL/hr"
$ws.Range("H23").Value = "5 mL/hr"
$ws.Range("I23").Value = "200-800-1000"
$ws.Range("J23").Value = "2 mins"
$ws.Range("K23").Value = "N"

$ws.Range("H23").Select()
